# CDS Input file updates
# Replace the "ParticipantsTab" Neo4j query (cell B2 on the "startup" sheet)
# with the updated version that sorts the collected sample ids via
# apoc.coll.sort before joining them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")
$ws.Activate()

$newQuery = "MATCH (p:participant)-->(s:study)`n" + `
    "OPTIONAL MATCH (samp:sample)-->(p)`n" + `
    "OPTIONAL MATCH (p)<--(diag:diagnosis)`n" + `
    "OPTIONAL MATCH (samp)<--(f:file)`n" + `
    "OPTIONAL MATCH (f)<--(g:genomic_info)`n" + `
    "WITH s, p, samp, f, g, diag`n" + `
    "WHERE f.file_type in ['TXT']`n" + `
    "with p`n" + `
    "OPTIONAL MATCH (p)-->(s:study)`n" + `
    "OPTIONAL MATCH (samp:sample)-->(p)`n" + `
    "WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp`n" + `
    "RETURN`n" + `
    "coalesce(p.participant_id,'') as ``Participant ID``,`n" + `
    "coalesce(s.study_name, '') as ``Study Name``,`n" + `
    "coalesce(s.phs_accession,'') as ``Accession``,`n" + `
    "coalesce(p.gender,'') as ``Gender``,`n" + `
    "coalesce(apoc.text.join(samp, ','), '') as ``Samples```n" + `
    "ORDER BY p.participant_id LIMIT 100"

$ws.Range("B2").Value = $newQuery

# Update the window's view state to match the author's session: the active
# cell moves from C5 to B5 and the sheet is scrolled down a few rows.
$ws.Range("B5").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
